$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("A RICCHIGIA SRL", "aricchigia@gmail.com"),
    @("A. DARBO AG", "katrin.widauer@darbo.at"),
    @("A. GANDOLA &amp; C. SPA", "gandola@gandola.it"),
    @("A. LOACKER SPA", "marketing@loacker.com"),
    @("A.D. SRL", "info@aiellobio.it"),
    @("A.O.C. SOCIETA' COOP. AGRICOLA", "info@calabriaaoc.it"),
    @("A.R. TARTUFI SRL", "commerciale@valnerinatartufi.com"),
    @("AB MAURI ITALY SPA", "italy.accounting@abmauri.com"),
    @("ACCADEMIA DEL PANE", "Not found/case 1"),
    @("ACCADEMIA GROUP SRL", "redazione@ristorazioneitalianamagazine.it"),
    @("ACCADEMIA OLEARIA SRL", "commerciale@accademiaolearia.com"),
    @("ACCUDIRE SRL", "info@accudire.eu"),
    @("ACEITUNAS SANMER OLIVES", "INFO@SANMER.COM"),
    @("ACEITUNAS TORREMAR", "Not found/case 2"),
    @("ACETAIA BORGO CASTELLO SRL", "avivici@acetaiaborgocastello.it")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $startRow + $i
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
}
